$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column price values are written as text (avoid Excel auto-numeric conversion)

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.883.35"
$ws.Range("E2").Value = "  +0.14%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.893.88"
$ws.Range("E3").Value = "  +0.01%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7727"
$ws.Range("E5").Value = "  -1.93%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "244.36"
$ws.Range("E6").Value = "  +0.42%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3129"
$ws.Range("E8").Value = "  -0.44%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.69"
$ws.Range("E9").Value = "  +1.70%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07224"
$ws.Range("E10").Value = "  -0.43%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08866"
$ws.Range("E11").Value = "  +9.57%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.003.36"
$ws.Range("E12").Value = "  +5.75%  "

# Row 13
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7718"
$ws.Range("E13").Value = "  +0.98%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.430"
$ws.Range("E14").Value = "  -1.80%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.44"
$ws.Range("E15").Value = "  +2.22%  "

# Row 16
$ws.Range("E16").Value = "  +1.07%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.886.30"
$ws.Range("E17").Value = "  +0.11%  "

# Row 18
$ws.Range("E18").Value = "  +0.37%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.63"
$ws.Range("E19").Value = "  +0.62%  "

# Row 20
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.282.03"
$ws.Range("E20").Value = "  +5.67%  "

# Row 21
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000007862"
$ws.Range("E21").Value = "  +1.10%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.135"
$ws.Range("E22").Value = "  +0.31%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.09%  "

# Row 24
$ws.Range("E24").Value = "  -0.08%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1602"
$ws.Range("E25").Value = "  -2.60%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.519"
$ws.Range("E26").Value = "  +1.44%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.35"
$ws.Range("E27").Value = "  -0.66%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.82"
$ws.Range("E28").Value = "  +0.63%  "

# Row 29
$ws.Range("E29").Value = "  -0.40%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.430"
$ws.Range("E30").Value = "  +2.15%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.544"
$ws.Range("E31").Value = "  -0.11%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.548"
$ws.Range("E32").Value = "  +1.85%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.120"
$ws.Range("E33").Value = "  +0.65%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05485"
$ws.Range("E34").Value = "  -0.78%  "

# Row 35
$ws.Range("E35").Value = "  -1.51%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7542"
$ws.Range("E36").Value = "  +2.15%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.004"
$ws.Range("E37").Value = "  +0.33%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.717"
$ws.Range("E38").Value = "  +3.31%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01960"
$ws.Range("E39").Value = "  +1.87%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.787"
$ws.Range("E40").Value = "  +0.35%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4512"
$ws.Range("E41").Value = "  +2.11%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.02"
$ws.Range("E42").Value = "  +0.18%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.092.87"
$ws.Range("E43").Value = "  -4.21%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.047"
$ws.Range("E44").Value = "  +2.97%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8555"
$ws.Range("E45").Value = "  +0.52%  "

# Row 46
$ws.Range("E46").Value = "  -0.06%  "

# Row 47
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.180.62"
$ws.Range("E47").Value = "  +6.22%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.15"
$ws.Range("E48").Value = "  -1.03%  "

# Row 49
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.888"
$ws.Range("E49").Value = "  +0.65%  "

# Row 50
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.609"
$ws.Range("E50").Value = "  +2.25%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.818"
$ws.Range("E51").Value = "  -1.69%  "
